# Refresh dashboard output: advance the "as of" date by one day and
# recompute the date-dependent pacing metrics on the Excel_vs_ML sheet.
#
# Columns (1-based):
#   C=Flight_Start_Date  D=Flight_End_Date      E=Today_Date
#   F=Last_Data_Till     G=Total_Budget          H=Spend_to_Date
#   I=Days_Elapsed       J=Days_Left             K=Expected_Spend_Till_Date
#   L=Pacing_%_vs_Ideal  M=Remaining_Budget      N=Ideal_Daily_Spend

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Excel_vs_ML")

$lastRow = 26

for ($row = 2; $row -le $lastRow; $row++) {
    $C = $ws.Cells.Item($row, 3).Value2
    $D = $ws.Cells.Item($row, 4).Value2
    $F = $ws.Cells.Item($row, 6).Value2
    $G = $ws.Cells.Item($row, 7).Value2
    $H = $ws.Cells.Item($row, 8).Value2

    # Advance "today" and "last data till" by one day.
    $newE = $F + 2
    $newF = $F + 1

    # Recompute the pacing metrics that depend on Last_Data_Till.
    $daysSinceStart = $newF - $C + 1
    $daysInFlight = $D - $C + 1
    $I = [Math]::Min($daysSinceStart, $daysInFlight)
    $J = [Math]::Max($D - $newF, 0)

    $K = [Math]::Round($G * $I / ($I + $J), 2)
    $M = $G - $H

    if ($K -ne 0) {
        $L = [Math]::Round(($H / $K) * 100, 2)
    } else {
        $L = 0
    }

    if ($J -ne 0) {
        $N = [Math]::Round($M / $J, 2)
    } else {
        $N = 0
    }

    $ws.Cells.Item($row, 5).Value = $newE
    $ws.Cells.Item($row, 6).Value = $newF
    $ws.Cells.Item($row, 9).Value = $I
    $ws.Cells.Item($row, 10).Value = $J
    $ws.Cells.Item($row, 11).Value = $K
    $ws.Cells.Item($row, 12).Value = $L
    $ws.Cells.Item($row, 14).Value = $N
}
